$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values for A2:A49 (LocalT2 re-solved results)
$values = @(
    85.009140767824505,
    85.877513711151735,
    86.060329067641689,
    89.213893967093242,
    89.762340036563074,
    89.625228519195616,
    80.484460694698356,
    81.901279707495434,
    81.261425959780624,
    80.347349177330898,
    79.296160877513714,
    81.946983546617915,
    80.255941499085921,
    80.438756855575861,
    82.541133455210229,
    80.0274223034735,
    81.581352833638022,
    87.020109689213896,
    90.447897623400365,
    90.85923217550274,
    90.813528336380259,
    83.638025594149909,
    87.340036563071294,
    86.791590493601461,
    83.912248628884825,
    83.455210237659955,
    83.957952468007306,
    82.266910420475327,
    81.53564899451554,
    81.307129798903105,
    90.356489945155388,
    93.327239488117002,
    92.961608775137123,
    87.568555758683729,
    91.910420475319938,
    82.861060329067641,
    88.71115173674589,
    85.191956124314444,
    84.232175502742237,
    81.307129798903105,
    81.855575868372938,
    81.398537477148082,
    81.672760511882998,
    84.689213893967093,
    86.106032906764156,
    80.530164533820837,
    82.358318098720289,
    81.352833638025587
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
